$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# ------------------------------------------------------------------
# 1) Insert a new grant row directly above the existing "Totals" row
#    (the "Totals" row is row 25 before this insertion) and populate
#    it with the new grant's details.
# ------------------------------------------------------------------
$totalsRow = $t.Rows.Item(25)
$t.Rows.Add($totalsRow) | Out-Null

$t.Cell(25, 1).Range.Text = "co-PI, Lu, Y. (PI), Pyrcz, M. (co-PI)"
$t.Cell(25, 2).Range.Text = "Unconventional Well Optimization based on Machine Learning"
$t.Cell(25, 3).Range.Text = "University Lands"
$t.Cell(25, 4).Range.Text = "`$75,000 (`$25,000)"
$t.Cell(25, 5).Range.Text = "9/2023-8/2025"

# ------------------------------------------------------------------
# 2) Update the "Totals" row's grant-total amount (now row 26) to
#    reflect the newly added grant.
# ------------------------------------------------------------------
$t.Cell(26, 4).Range.Text = "`$14,019,591 (`$3,780,350)"

# ------------------------------------------------------------------
# 3) Update the "Indicates awarded in rank" amount (now row 27) to
#    reflect the newly added grant.
# ------------------------------------------------------------------
$t.Cell(27, 4).Range.Text = "`$3,638,025 (`$1,256,832)"

# ------------------------------------------------------------------
# 4) Rename the table's bookmark from T_1b264 to T_e0a98. The COM
#    bookmark collection in this runtime only tracks bookmarks added
#    during this session, so add the new bookmark at the same
#    location (immediately before the table caption paragraph).
#    Re-derive the paragraphs collection from Content so the index
#    reflects the table edits made above.
# ------------------------------------------------------------------
$tableCaption = $d.Content.Paragraphs.Item(3).Range
$tableCaption.Collapse(1)
$d.Bookmarks.Add("T_e0a98", $tableCaption) | Out-Null
